$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 9).Value = "b"
$ws.Cells.Item(4, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(15, 9).Value = "%"
$ws.Cells.Item(15, 10).Value = "Uninterpretable"
$ws.Cells.Item(16, 9).Value = "sd"
$ws.Cells.Item(16, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(22, 9).Value = "b"
$ws.Cells.Item(22, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(38, 9).Value = "sd"
$ws.Cells.Item(38, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(47, 9).Value = "ba"
$ws.Cells.Item(47, 10).Value = "Appreciation"
$ws.Cells.Item(60, 9).Value = "sv"
$ws.Cells.Item(60, 10).Value = "Statement-opinion"
$ws.Cells.Item(62, 9).Value = "b"
$ws.Cells.Item(62, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(63, 9).Value = "%"
$ws.Cells.Item(63, 10).Value = "Uninterpretable"
$ws.Cells.Item(71, 9).Value = "sv"
$ws.Cells.Item(71, 10).Value = "Statement-opinion"
$ws.Cells.Item(75, 9).Value = "%"
$ws.Cells.Item(75, 10).Value = "Uninterpretable"
$ws.Cells.Item(91, 9).Value = "ba"
$ws.Cells.Item(91, 10).Value = "Appreciation"
$ws.Cells.Item(92, 9).Value = "sv"
$ws.Cells.Item(92, 10).Value = "Statement-opinion"
$ws.Cells.Item(99, 9).Value = "qy"
$ws.Cells.Item(99, 10).Value = "Yes-No-Question"
$ws.Cells.Item(101, 9).Value = "b"
$ws.Cells.Item(101, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(146, 9).Value = "%"
$ws.Cells.Item(146, 10).Value = "Uninterpretable"
$ws.Cells.Item(147, 9).Value = "%"
$ws.Cells.Item(147, 10).Value = "Uninterpretable"
$ws.Cells.Item(160, 9).Value = "sv"
$ws.Cells.Item(160, 10).Value = "Statement-opinion"
$ws.Cells.Item(178, 9).Value = "sd"
$ws.Cells.Item(178, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(206, 9).Value = "b"
$ws.Cells.Item(206, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(223, 9).Value = "sv"
$ws.Cells.Item(223, 10).Value = "Statement-opinion"
$ws.Cells.Item(224, 9).Value = "sd"
$ws.Cells.Item(224, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(232, 9).Value = "sv"
$ws.Cells.Item(232, 10).Value = "Statement-opinion"
$ws.Cells.Item(248, 9).Value = "sv"
$ws.Cells.Item(248, 10).Value = "Statement-opinion"
$ws.Cells.Item(254, 9).Value = "sv"
$ws.Cells.Item(254, 10).Value = "Statement-opinion"
